$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R4RPageLoad")

# Remove the "R4R Home" (row 2) and "R4R Results" (row 3) load events -
# these were being logged incorrectly. The remaining two "R4R Detail" rows
# shift up to become rows 2 and 3.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# Update the selection to sit just below the remaining data.
$ws.Range("A4").Select()
